# Rewrite the Angpt2-Tek NATMI edge table (Sheet1, A1:T17) with the
# refreshed TPM-derived values. The pipeline re-run added the
# "Resolving-Mac" sending cluster (previously missing), so every
# sending/target cluster pairing is regenerated below (header + 16 data
# rows spanning the 4x4 ECs/FAPs/MuSCs/Resolving-Mac cluster cross).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @('Sending cluster', 'Ligand symbol', 'Receptor symbol', 'Target cluster', 'Ligand-expressing cells', 'Ligand detection rate', 'Ligand average expression value', 'Ligand total expression value', 'Ligand derived specificity of average expression value', 'Ligand derived specificity of total expression value', 'Receptor-expressing cells', 'Receptor detection rate', 'Receptor average expression value', 'Receptor total expression value', 'Receptor derived specificity of average expression value', 'Receptor derived specificity of total expression value', 'Edge average expression weight', 'Edge total expression weight', 'Edge average expression derived specificity', 'Edge total expression derived specificity'),
        @('ECs', 'Angpt2', 'Tek', 'ECs', 3, 1, 9.450759333333334, 28.352278, 0.6100581096599255, 0.6100581096599255, 3, 1, 56.16842399999999, 168.505272, 0.8021403195141494, 0.8021403195141494, 530.834257356624, 4777.508316209616, 0.4893522070048106, 0.4893522070048106),
        @('ECs', 'Angpt2', 'Tek', 'FAPs', 3, 1, 9.450759333333334, 28.352278, 0.6100581096599255, 0.6100581096599255, 3, 1, 10.582537, 31.747611, 0.1511290331103168, 0.1511290331103167, 100.0130103230953, 900.117092907858, 0.09219749225401212, 0.09219749225401211),
        @('ECs', 'Angpt2', 'Tek', 'MuSCs', 3, 1, 9.450759333333334, 28.352278, 0.6100581096599255, 0.6100581096599255, 3, 1, 3.262797333333333, 9.788392, 0.04659595390231912, 0.04659595390231912, 30.83591235077511, 277.523211156976, 0.02842623955544983, 0.02842623955544983),
        @('ECs', 'Angpt2', 'Tek', 'Resolving-Mac', 3, 1, 9.450759333333334, 28.352278, 0.6100581096599255, 0.6100581096599255, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.0891364117788889, 0.80222770601, 0.00008217084565283581, 0.00008217084565283584),
        @('FAPs', 'Angpt2', 'Tek', 'ECs', 3, 1, 1.769696333333333, 5.309089, 0.1142360694740755, 0.1142360694740756, 3, 1, 56.16842399999999, 168.505272, 0.8021403195141494, 0.8021403195141494, 99.40105400191199, 894.6094860172079, 0.09163335726797553, 0.09163335726797554),
        @('FAPs', 'Angpt2', 'Tek', 'FAPs', 3, 1, 1.769696333333333, 5.309089, 0.1142360694740755, 0.1142360694740756, 3, 1, 10.582537, 31.747611, 0.1511290331103168, 0.1511290331103167, 18.72787692626433, 168.550892336379, 0.01726438672594001, 0.01726438672594001),
        @('FAPs', 'Angpt2', 'Tek', 'MuSCs', 3, 1, 1.769696333333333, 5.309089, 0.1142360694740755, 0.1142360694740756, 3, 1, 3.262797333333333, 9.788392, 0.04659595390231912, 0.04659595390231912, 5.774160477209778, 51.967444294888, 0.005322938627196148, 0.005322938627196149),
        @('FAPs', 'Angpt2', 'Tek', 'Resolving-Mac', 3, 1, 1.769696333333333, 5.309089, 0.1142360694740755, 0.1142360694740756, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.01669118591722222, 0.150220673255, 0.00001538685296384891, 0.00001538685296384892),
        @('MuSCs', 'Angpt2', 'Tek', 'ECs', 3, 1, 3.398621333333333, 10.195864, 0.2193851766757396, 0.2193851766757396, 3, 1, 56.16842399999999, 168.505272, 0.8021403195141494, 0.8021403195141494, 190.895204066112, 1718.056836595008, 0.1759776957153459, 0.1759776957153459),
        @('MuSCs', 'Angpt2', 'Tek', 'FAPs', 3, 1, 3.398621333333333, 10.195864, 0.2193851766757396, 0.2193851766757396, 3, 1, 10.582537, 31.747611, 0.1511290331103168, 0.1511290331103167, 35.96603600898933, 323.694324080904, 0.03315546962974054, 0.03315546962974054),
        @('MuSCs', 'Angpt2', 'Tek', 'MuSCs', 3, 1, 3.398621333333333, 10.195864, 0.2193851766757396, 0.2193851766757396, 3, 1, 3.262797333333333, 9.788392, 0.04659595390231912, 0.04659595390231912, 11.08901262340978, 99.80111361068801, 0.0102224615792349, 0.0102224615792349),
        @('MuSCs', 'Angpt2', 'Tek', 'Resolving-Mac', 3, 1, 3.398621333333333, 10.195864, 0.2193851766757396, 0.2193851766757396, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.03205466354222222, 0.28849197188, 0.00002954975141825659, 0.0000295497514182566),
        @('Resolving-Mac', 'Angpt2', 'Tek', 'ECs', 3, 1, 0.8724953333333333, 2.617486, 0.0563206441902594, 0.05632064419025941, 3, 1, 56.16842399999999, 168.505272, 0.8021403195141494, 0.8021403195141494, 49.00668782068799, 441.0601903861919, 0.0451770595260174, 0.04517705952601741),
        @('Resolving-Mac', 'Angpt2', 'Tek', 'FAPs', 3, 1, 0.8724953333333333, 2.617486, 0.0563206441902594, 0.05632064419025941, 3, 1, 10.582537, 31.747611, 0.1511290331103168, 0.1511290331103167, 9.233214147327333, 83.09892732594599, 0.008511684500624082, 0.008511684500624081),
        @('Resolving-Mac', 'Angpt2', 'Tek', 'MuSCs', 3, 1, 0.8724953333333333, 2.617486, 0.0563206441902594, 0.05632064419025941, 3, 1, 3.262797333333333, 9.788392, 0.04659595390231912, 0.04659595390231912, 2.846775446945778, 25.620979022512, 0.002624314140438244, 0.002624314140438245),
        @('Resolving-Mac', 'Angpt2', 'Tek', 'Resolving-Mac', 3, 1, 0.8724953333333333, 2.617486, 0.0563206441902594, 0.05632064419025941, 1, 0.3333333333333333, 0.009431666666666666, 0.028295, 0.0001346934732146117, 0.0001346934732146117, 0.008229085152222222, 0.07406176637, 0.00000758602317967038, 0.000007586023179670383)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

Write-Output "Wrote $($data.Length) rows (A1:T$($data.Length)) to $($ws.Name)"
